$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sneha's image reference from .png to .jpg
$ws.Range("B15").Value = "meme_images/sneha.jpg"

# View state changed: sheet scrolled so row 6 is the top visible row,
# and the active selection moved to F17.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("F17").Select()
